# feat: add 2022-Q4 data
#
# 1. Insert a new "2022-Q4" worksheet right after "总计" and before "2022-Q3",
#    populated with 8 rows of fund-holding data (cloning the formatting of
#    the existing "2022-Q3" sheet so headers/borders/fonts match).
# 2. Insert a new row into the "总计" (summary) sheet for the "2022-Q4"
#    quarter, pushing the existing quarters' figures down by one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Part 1: "总计" summary sheet - insert new row for 2022-Q4
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Push existing data rows (2..7) down to (3..8).
$summary.Rows.Item(2).Insert()

# Copy the (now shifted) row 3 formatting into the new row 2 so the new
# row matches the existing look (borders/fonts/number formats).
$summary.Range("A3:D3").Copy($summary.Range("A2:D2"))

# Populate the new summary row for 2022-Q4.
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 8
$summary.Range("D2").Value = 2.32

# Column A is a plain sequential row index (0,1,2,...) - it is NOT shifted
# with the other data, it always just counts up. Re-assert it for every
# data row so it reads 0..6 top to bottom.
$summary.Range("A2").Value = 0
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4
$summary.Range("A7").Value = 5
$summary.Range("A8").Value = 6

# ---------------------------------------------------------------------
# Part 2: brand-new "2022-Q4" worksheet with fund holdings
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"
$wb.Worksheets.Item("2022-Q4").Move($wb.Worksheets.Item("2022-Q3"))

# The template sheet only has one data row (row 2); we need 8 (rows 2-9).
# Insert 7 more rows right below the existing data row, then stamp the
# row-2 formatting onto them so every row matches.
for ($i = 0; $i -lt 7; $i++) {
    $q4.Rows.Item(3).Insert()
}
$q4.Range("A2:H2").Copy($q4.Range("A3:H9"))

# --- row 2 ---
$q4.Range("A2").Value = 0
$q4.Range("B2").NumberFormat = "@"
$q4.Range("B2").Value = "100060"
$q4.Range("C2").Value = "富国高新技术产业混合"
$q4.Range("D2").NumberFormat = "@"
$q4.Range("D2").Value = "14.62"
$q4.Range("E2").NumberFormat = "@"
$q4.Range("E2").Value = "93.20"
$q4.Range("F2").NumberFormat = "@"
$q4.Range("F2").Value = "5.07"
$q4.Range("G2").NumberFormat = "@"
$q4.Range("G2").Value = "0.7412"
$q4.Range("H2").Value = 2

# --- row 3 ---
$q4.Range("A3").Value = 1
$q4.Range("B3").NumberFormat = "@"
$q4.Range("B3").Value = "009863"
$q4.Range("C3").Value = "富国创新趋势股票"
$q4.Range("D3").NumberFormat = "@"
$q4.Range("D3").Value = "30.80"
$q4.Range("E3").NumberFormat = "@"
$q4.Range("E3").Value = "92.72"
$q4.Range("F3").NumberFormat = "@"
$q4.Range("F3").Value = "2.18"
$q4.Range("G3").NumberFormat = "@"
$q4.Range("G3").Value = "0.6714"
$q4.Range("H3").Value = 8

# --- row 4 ---
$q4.Range("A4").Value = 2
$q4.Range("B4").NumberFormat = "@"
$q4.Range("B4").Value = "007345"
$q4.Range("C4").Value = "富国科技创新灵活配置混合"
$q4.Range("D4").NumberFormat = "@"
$q4.Range("D4").Value = "8.84"
$q4.Range("E4").NumberFormat = "@"
$q4.Range("E4").Value = "94.61"
$q4.Range("F4").NumberFormat = "@"
$q4.Range("F4").Value = "5.00"
$q4.Range("G4").NumberFormat = "@"
$q4.Range("G4").Value = "0.4420"
$q4.Range("H4").Value = 1

# --- row 5 ---
$q4.Range("A5").Value = 3
$q4.Range("B5").NumberFormat = "@"
$q4.Range("B5").Value = "519021"
$q4.Range("C5").Value = "国泰金鼎价值混合"
$q4.Range("D5").NumberFormat = "@"
$q4.Range("D5").Value = "5.86"
$q4.Range("E5").NumberFormat = "@"
$q4.Range("E5").Value = "86.87"
$q4.Range("F5").NumberFormat = "@"
$q4.Range("F5").Value = "3.71"
$q4.Range("G5").NumberFormat = "@"
$q4.Range("G5").Value = "0.2174"
$q4.Range("H5").Value = 8

# --- row 6 ---
$q4.Range("A6").Value = 4
$q4.Range("B6").NumberFormat = "@"
$q4.Range("B6").Value = "016021"
$q4.Range("C6").Value = "华安优嘉精选混合A"
$q4.Range("D6").NumberFormat = "@"
$q4.Range("D6").Value = "6.09"
$q4.Range("E6").NumberFormat = "@"
$q4.Range("E6").Value = "61.72"
$q4.Range("F6").NumberFormat = "@"
$q4.Range("F6").Value = "1.86"
$q4.Range("G6").NumberFormat = "@"
$q4.Range("G6").Value = "0.1133"
$q4.Range("H6").Value = 10

# --- row 7 ---
$q4.Range("A7").Value = 5
$q4.Range("B7").NumberFormat = "@"
$q4.Range("B7").Value = "016022"
$q4.Range("C7").Value = "华安优嘉精选混合C"
$q4.Range("D7").NumberFormat = "@"
$q4.Range("D7").Value = "5.26"
$q4.Range("E7").NumberFormat = "@"
$q4.Range("E7").Value = "61.72"
$q4.Range("F7").NumberFormat = "@"
$q4.Range("F7").Value = "1.86"
$q4.Range("G7").NumberFormat = "@"
$q4.Range("G7").Value = "0.0978"
$q4.Range("H7").Value = 10

# --- row 8 ---
$q4.Range("A8").Value = 6
$q4.Range("B8").NumberFormat = "@"
$q4.Range("B8").Value = "001626"
$q4.Range("C8").Value = "国泰央企改革股票"
$q4.Range("D8").NumberFormat = "@"
$q4.Range("D8").Value = "0.54"
$q4.Range("E8").NumberFormat = "@"
$q4.Range("E8").Value = "86.52"
$q4.Range("F8").NumberFormat = "@"
$q4.Range("F8").Value = "3.71"
$q4.Range("G8").NumberFormat = "@"
$q4.Range("G8").Value = "0.0200"
$q4.Range("H8").Value = 8

# --- row 9 ---
$q4.Range("A9").Value = 7
$q4.Range("B9").NumberFormat = "@"
$q4.Range("B9").Value = "006836"
$q4.Range("C9").Value = "永赢惠泽一年定期开放灵活配置混合"
$q4.Range("D9").NumberFormat = "@"
$q4.Range("D9").Value = "3.73"
$q4.Range("E9").NumberFormat = "@"
$q4.Range("E9").Value = "48.14"
$q4.Range("F9").NumberFormat = "@"
$q4.Range("F9").Value = "0.41"
$q4.Range("G9").NumberFormat = "@"
$q4.Range("G9").Value = "0.0153"
$q4.Range("H9").Value = 7

# Clear the temporary text-number-format so the cells end up with the
# same (default) style as the rest of the sheet, matching the source.
$q4.Range("B2:G9").ClearFormats()

Write-Host "2022-Q4 sheet + summary row added"
